# Backlog_19.xlsx: mark the "ITI" week-19 items that were closed out, and the
# two matching "SPN" items, as "Resolvido" instead of "Pendente" (Status
# column = J). Also refresh the saved view state (active sheet/selection)
# to reflect where the user ended up working.

$wb = $excel.ActiveWorkbook

# --- ITI sheet: rows 2,3,5,7,10,11,12,15 move from "Pendente" to "Resolvido" ---
$wsITI = $wb.Worksheets.Item("ITI")
$itiResolvedRows = @(2, 3, 5, 7, 10, 11, 12, 15)
foreach ($row in $itiResolvedRows) {
    $wsITI.Range("J$row").Value = "Resolvido"
}

# --- SPN sheet: rows 2 and 3 move from "Pendente" to "Resolvido" ---
$wsSPN = $wb.Worksheets.Item("SPN")
$spnResolvedRows = @(2, 3)
foreach ($row in $spnResolvedRows) {
    $wsSPN.Range("J$row").Value = "Resolvido"
}

# --- Restore the saved selection on ITI (no longer the active tab) ---
$wsITI.Activate()
$wsITI.Range("C17").Select()

# --- SPN becomes the active tab, with C19:C20 selected ---
$wsSPN.Activate()
$wsSPN.Range("C19:C20").Select()
